$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title (row 1) -----------------------------------------------------
# A1 holds a rich-text string: "Table 1" (bold) + ": Risk factor definitions. ..."
# Only the non-bold tail changes; re-apply Bold to the first run so the
# rich-text split survives the round trip.
$titleCell = $ws.Range("A1")
$fullTitle = $titleCell.Value2
$tailLen = $fullTitle.Length - 7
$tailChars = $titleCell.Characters(8, $tailLen)
$tailChars.Text = ": Risk factor definitions with categories are marked with (*)."

$headChars = $titleCell.Characters(1, 7)
$headChars.Font.Bold = $true
$headChars.Font.Name = "Calibri"
$headChars.Font.Size = 11
$tailChars2 = $titleCell.Characters(8, $titleCell.Value2.Length - 7)
$tailChars2.Font.Name = "Calibri"
$tailChars2.Font.Size = 11

# --- BMI category label (row 7) ----------------------------------------
$ws.Range("A7").Value = "Body Mass Index (BMI)"

# --- Physical activity definitions (rows 11-13) -------------------------
$ws.Range("C11").Value = "0 to <1.5 MET-hours/day"
$ws.Range("C12").Value = "1.5 to <3 MET-hours/day"
$ws.Range("C13").Value = "≥ 3 MET-hours/day"

# --- Wrap text for the BMI category's merged block (A7:A10) -------------
# The longer "Body Mass Index (BMI)" label now wraps inside its cell.
$ws.Range("A7:A10").WrapText = $true
$ws.Rows.Item(7).RowHeight = 15.75

# --- Scroll the view down a bit (cosmetic) -------------------------------
try {
    $excel.ActiveWindow.ScrollRow = 10
} catch {
}
